$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (Leetcode Question No. / Question / Language)
$ws.Range("A27").Value = "GFG"
$ws.Range("B27").Value = "Intersection Point in Y Shaped Linked Lists"
$ws.Range("C27").Value = "Java"

# Match the author's final view state: scrolled down with E32 selected
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("E32").Select()
